$wb = $excel.ActiveWorkbook

$wsMultiLogin  = $wb.Worksheets.Item("MultiLogin")
$wsMultiLogin2 = $wb.Worksheets.Item("MultiLogin2")

# --- MultiLogin (sheet2.xml): shift the account-type / username rows down
#     (a "Global Admin" row is now the first data row) and update the
#     selection. Cell A2's format moves to the "Global Admin" look used on
#     MultiLogin2's A2, so pull the formatting across before overwriting
#     the values.
$wsMultiLogin2.Range("A2").Copy()
$wsMultiLogin.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsMultiLogin.Range("A2").Value = "Global Admin"
$wsMultiLogin.Range("B2").Value = "changepondtest16@yahoo.com"
$wsMultiLogin.Range("A3").Value = "Shipper Admin"
$wsMultiLogin.Range("B3").Value = "karthikeyan.s@changepond.com"
$wsMultiLogin.Range("A4").Value = "Shipper User"
$wsMultiLogin.Range("B4").Value = "karthirko29@gmail.com"

# --- MultiLogin2 (sheet3.xml): selection moves from A4:C4 to A2:C3.
#     Select it on that sheet without leaving it as the active tab, then
#     restore MultiLogin's own selection/active-tab state.
$wsMultiLogin2.Range("A2:C3").Select()

$wsMultiLogin.Activate()
$wsMultiLogin.Range("A6").Select()
